{"js": "// Replace the worksheet date and the 25 multiplication problems with the\n// new values from the updated OOXML. Every source string below is unique\n// within the document, so a plain text search+replace is unambiguous.\nconst replacements = [\n  [\"2025-01-20 Monday\", \"2025-01-21 Tuesday\"],\n  [\"191\u00d75=\", \"507\u00d75=\"],\n  [\"976\u00d72=\", \"674\u00d77=\"],\n  [\"137\u00d75=\", \"763\u00d79=\"],\n  [\"872\u00d75=\", \"444\u00d77=\"],\n  [\"296\u00d74=\", \"569\u00d74=\"],\n  [\"612\u00d73=\", \"704\u00d77=\"],\n  [\"857\u00d74=\", \"603\u00d79=\"],\n  [\"455\u00d79=\", \"542\u00d76=\"],\n  [\"820\u00d73=\", \"594\u00d72=\"],\n  [\"670\u00d75=\", \"146\u00d78=\"],\n  [\"559\u00d73=\", \"933\u00d78=\"],\n  [\"473\u00d76=\", \"149\u00d74=\"],\n  [\"192\u00d74=\", \"291\u00d79=\"],\n  [\"314\u00d76=\", \"402\u00d75=\"],\n  [\"319\u00d77=\", \"127\u00d79=\"],\n  [\"508\u00d77=\", \"246\u00d72=\"],\n  [\"888\u00d79=\", \"679\u00d75=\"],\n  [\"961\u00d78=\", \"822\u00d73=\"],\n  [\"554\u00d77=\", \"145\u00d77=\"],\n  [\"765\u00d72=\", \"487\u00d72=\"],\n  [\"783\u00d76=\", \"996\u00d73=\"],\n  [\"530\u00d76=\", \"108\u00d77=\"],\n  [\"117\u00d73=\", \"786\u00d78=\"],\n  [\"819\u00d78=\", \"551\u00d74=\"],\n  [\"738\u00d79=\", \"966\u00d76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the worksheet date and the 25 multiplication problems with the\n# new values from the updated OOXML. Every source string is unique within\n# the document, so Find/Replace (ReplaceAll) is unambiguous for each pair.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2025-01-20 Monday\", \"2025-01-21 Tuesday\"),\n  @(\"191\u00d75=\", \"507\u00d75=\"),\n  @(\"976\u00d72=\", \"674\u00d77=\"),\n  @(\"137\u00d75=\", \"763\u00d79=\"),\n  @(\"872\u00d75=\", \"444\u00d77=\"),\n  @(\"296\u00d74=\", \"569\u00d74=\"),\n  @(\"612\u00d73=\", \"704\u00d77=\"),\n  @(\"857\u00d74=\", \"603\u00d79=\"),\n  @(\"455\u00d79=\", \"542\u00d76=\"),\n  @(\"820\u00d73=\", \"594\u00d72=\"),\n  @(\"670\u00d75=\", \"146\u00d78=\"),\n  @(\"559\u00d73=\", \"933\u00d78=\"),\n  @(\"473\u00d76=\", \"149\u00d74=\"),\n  @(\"192\u00d74=\", \"291\u00d79=\"),\n  @(\"314\u00d76=\", \"402\u00d75=\"),\n  @(\"319\u00d77=\", \"127\u00d79=\"),\n  @(\"508\u00d77=\", \"246\u00d72=\"),\n  @(\"888\u00d79=\", \"679\u00d75=\"),\n  @(\"961\u00d78=\", \"822\u00d73=\"),\n  @(\"554\u00d77=\", \"145\u00d77=\"),\n  @(\"765\u00d72=\", \"487\u00d72=\"),\n  @(\"783\u00d76=\", \"996\u00d73=\"),\n  @(\"530\u00d76=\", \"108\u00d77=\"),\n  @(\"117\u00d73=\", \"786\u00d78=\"),\n  @(\"819\u00d78=\", \"551\u00d74=\"),\n  @(\"738\u00d79=\", \"966\u00d76=\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $newText\n  $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
